$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Mina Panulcillo"
$ws.Range("C3").Value = 5

$ws.Range("E5").Select()
